$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accountability rows (8-11) - camera / RC gear purchased by Juan.
# NOTE: column B is set before column A in row 8 so that the new shared
# strings land in the same order as the target workbook
# ("Canon Elph 330 S" = index 20, "Canon" = index 21, ...).
$ws.Range("B8").Value = "Canon Elph 330 S"
$ws.Range("A8").Value = "Canon"
$ws.Range("C8").Value = 167.13
$ws.Range("D8").Value = "Juan"

$ws.Range("A9").Value = "Hobby king"
$ws.Range("B9").Value = "RC 6 channels"
$ws.Range("C9").Value = 71.57
$ws.Range("D9").Value = "Juan"

$ws.Range("A10").Value = "StockRC"
$ws.Range("B10").Value = "Battery charger"
$ws.Range("C10").Value = 68.62
$ws.Range("D10").Value = "Juan"

$ws.Range("A11").Value = "StockRC"
$ws.Range("B11").Value = "Servo wire x5"
$ws.Range("C11").Value = 14.54
$ws.Range("D11").Value = "Juan"

# Update the current selection to match the author's final cursor position.
[void]$ws.Range("D13").Select()
